# Add more Modbus-TCP registers for the generator start/stop feature.
# Inserts three new rows (State, Error, NoGeneratorAtAcIn) right after the
# existing "/Generator0/Runtime" (uint32) row on the "Field list" sheet,
# pushing the pre-existing com.victronenergy.meteo rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Field list")

# Insert 3 blank rows before the current row 351 (the first meteo row).
$ws.Range("A351:A353").EntireRow.Insert()

# Match the compact row height used throughout this table (13.8pt), same as
# the row immediately above and the rows that just got pushed down.
$ws.Rows.Item(351).RowHeight = 13.8
$ws.Rows.Item(352).RowHeight = 13.8
$ws.Rows.Item(353).RowHeight = 13.8

# --- Row 351: Generator start/stop state ---------------------------------
$ws.Cells.Item(351, 1).Value = "com.victronenergy.generator"
$ws.Cells.Item(351, 2).Value = "Generator start/stop state"
$ws.Cells.Item(351, 3).Value = 3506
$ws.Cells.Item(351, 4).Value = "uint16"
$ws.Cells.Item(351, 5).Value = 1
$ws.Cells.Item(351, 6).Value = "0 to 65535"
$ws.Cells.Item(351, 7).Value = "/Generator0/State"
$ws.Cells.Item(351, 8).Value = "no"
$ws.Cells.Item(351, 9).Value = "0=Stopped;1=Running;10=Error"

# --- Row 352: Generator remote error --------------------------------------
$ws.Cells.Item(352, 1).Value = "com.victronenergy.generator"
$ws.Cells.Item(352, 2).Value = "Generator remote error"
$ws.Cells.Item(352, 3).Value = 3507
$ws.Cells.Item(352, 4).Value = "uint16"
$ws.Cells.Item(352, 5).Value = 1
$ws.Cells.Item(352, 6).Value = "0 to 65535"
$ws.Cells.Item(352, 7).Value = "/Generator0/Error"
$ws.Cells.Item(352, 8).Value = "no"
$ws.Cells.Item(352, 9).Value = "0=No Error;1=Remote disabled;2=Remote fault"
$ws.Cells.Item(352, 10).Value = "Only used for FisherPanda gensets"

# --- Row 353: Generator not detected at AC input alarm --------------------
$ws.Cells.Item(353, 1).Value = "com.victronenergy.generator"
$ws.Cells.Item(353, 2).Value = "Generator not detected at AC input alarm"
$ws.Cells.Item(353, 3).Value = 3508
$ws.Cells.Item(353, 4).Value = "uint16"
$ws.Cells.Item(353, 5).Value = 1
$ws.Cells.Item(353, 6).Value = "0 to 65535"
$ws.Cells.Item(353, 7).Value = "/Generator0/Alarms/NoGeneratorAtAcIn"
$ws.Cells.Item(353, 8).Value = "no"
$ws.Cells.Item(353, 9).Value = "0=No alarm;2=Alarm"

# The pre-existing com.victronenergy.meteo rows (now pushed down to
# 354-357) previously had no explicit "writable" (column H) value; the
# author filled that column in for them as part of this same edit.
$ws.Cells.Item(354, 8).Value = "no"
$ws.Cells.Item(355, 8).Value = "no"
$ws.Cells.Item(356, 8).Value = "no"
$ws.Cells.Item(357, 8).Value = "no"

# Keep the view in sync with where the author left off editing (frozen pane
# scrolled down near the newly-added rows, selection parked on the last new
# generator row before the pre-existing meteo block).
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A329").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A353").Select()
